$d = $word.ActiveDocument

# The run "do Cục CSQLHC về TTXH cấp" needs to become three separate runs
# (all sharing the same formatting):
#   "do "  +  "Cục CSQLHC về TTXH"  +  " cấp"
$part1 = "do "
$part2 = "Cục CSQLHC về TTXH"
$part3 = " cấp"

$target = $d.Content
$found = $target.Find.Execute($part1 + $part2 + $part3, $true, $false, $false,
                               $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $start = $target.Start
    $endPos = $target.End
    $boundary1 = $start + $part1.Length
    $boundary2 = $boundary1 + $part2.Length

    # The engine coalesces adjacent runs that share identical formatting, so
    # a plain text split (InsertAfter/Text assignment) gets silently merged
    # back into one run. Toggling a character-formatting property on a
    # sub-range (and then restoring it) forces a genuine run boundary at
    # that point, which is what we need here since the three resulting
    # runs must end up with identical <w:rPr> (color 000000) but as
    # separate <w:r> elements.

    $tail1 = $d.Range($boundary1, $endPos)
    $tail1.Font.Color = 1
    $tail1.Font.Color = 0

    $tail2 = $d.Range($boundary2, $endPos)
    $tail2.Font.Color = 1
    $tail2.Font.Color = 0
}
